# Update attendance ("F") and price ("G") figures with freshly generated
# data, matching the gh-pages output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "展览" (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 7634
$ws1.Range("F3").Value  = 3596
$ws1.Range("F5").Value  = 3911
$ws1.Range("F7").Value  = 100
$ws1.Range("F8").Value  = 89
$ws1.Range("F9").Value  = 123
$ws1.Range("F10").Value = 188
$ws1.Range("F11").Value = 528
$ws1.Range("F13").Value = 171
$ws1.Range("F15").Value = 4
$ws1.Range("F16").Value = 11
$ws1.Range("F18").Value = 363
$ws1.Range("F19").Value = 4279
$ws1.Range("F20").Value = 4279
$ws1.Range("F23").Value = 1044
$ws1.Range("F25").Value = 2647
$ws1.Range("G25").Value = 72
$ws1.Range("F27").Value = 115
$ws1.Range("F28").Value = 3127
$ws1.Range("F29").Value = 2421
$ws1.Range("F30").Value = 78
$ws1.Range("F34").Value = 124
$ws1.Range("F36").Value = 47
$ws1.Range("F38").Value = 4558
$ws1.Range("F39").Value = 551
$ws1.Range("F42").Value = 942
$ws1.Range("F43").Value = 882
$ws1.Range("F44").Value = 279
$ws1.Range("F46").Value = 1720
$ws1.Range("F47").Value = 271
$ws1.Range("F50").Value = 744

# ---------------------------------------------------------------------------
# Sheet 2: "演出" (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F24").Value = 640

# ---------------------------------------------------------------------------
# Sheet 3: "本地生活" (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")

$ws3.Range("F2").Value = 174

# ---------------------------------------------------------------------------
# Sheet 4: "全部类型" (All types)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 174
$ws4.Range("F4").Value  = 7634
$ws4.Range("F5").Value  = 3596
$ws4.Range("F6").Value  = 3911
$ws4.Range("F7").Value  = 100
$ws4.Range("F8").Value  = 89
$ws4.Range("F9").Value  = 123
$ws4.Range("F11").Value = 188
$ws4.Range("F12").Value = 528
$ws4.Range("F14").Value = 171
$ws4.Range("F15").Value = 11
$ws4.Range("F17").Value = 363
$ws4.Range("F18").Value = 4279
$ws4.Range("F19").Value = 4279
$ws4.Range("F25").Value = 2647
$ws4.Range("G25").Value = 72
$ws4.Range("F27").Value = 115
$ws4.Range("F28").Value = 2421
$ws4.Range("F29").Value = 78
$ws4.Range("F31").Value = 124
$ws4.Range("F33").Value = 47
$ws4.Range("F37").Value = 4558
$ws4.Range("F39").Value = 551
$ws4.Range("F42").Value = 942
$ws4.Range("F43").Value = 882
$ws4.Range("F44").Value = 279
$ws4.Range("F46").Value = 1720
$ws4.Range("F47").Value = 271
$ws4.Range("F50").Value = 744
